# Update countries & provincias Spain
# Applies updated case statistics for several countries, and re-ranks
# Jamaica / Malta (Malta overtakes Jamaica so it moves one row up the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Brasil (row 11)
$ws.Range("B11").Value = 163427
$ws.Range("C11").Value = 728
$ws.Range("E11").Value = 87302
$ws.Range("G11").Value = 45
$ws.Range("H11").Value = 11168

# India (row 16)
$ws.Range("B16").Value = 67701
$ws.Range("C16").Value = 540
$ws.Range("D16").Value = 21147
$ws.Range("E16").Value = 44339

# Suiza (row 23)
$ws.Range("B23").Value = 30344
$ws.Range("C23").Value = 39
$ws.Range("E23").Value = 1910
$ws.Range("F23").Value = 89
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 1834

# Uzbekistan (row 75)
$ws.Range("D75").Value = 1982
$ws.Range("E75").Value = 461

# Bosnia y Herzegovina (row 78)
$ws.Range("B78").Value = 2141
$ws.Range("C78").Value = 24
$ws.Range("D78").Value = 1114
$ws.Range("E78").Value = 914
$ws.Range("G78").Value = 6
$ws.Range("H78").Value = 113

# Libano (row 105)
$ws.Range("B105").Value = 859
$ws.Range("C105").Value = 14
$ws.Range("E105").Value = 599

# Jamaica / Malta swap places (row 123 becomes Malta, row 124 becomes Jamaica)
# and Malta's figures are refreshed with newer totals.
$ws.Range("A123").Value = "Malta"
$ws.Range("B123").Value = 503
$ws.Range("C123").Value = 7
$ws.Range("D123").Value = 434
$ws.Range("E123").Value = 64
$ws.Range("H123").Value = 5

$ws.Range("A124").Value = "Jamaica"
$ws.Range("B124").Value = 502
$ws.Range("C124").Value = 4
$ws.Range("D124").Value = 90
$ws.Range("E124").Value = 403
$ws.Range("H124").Value = 9

# Nepal (row 158)
$ws.Range("D158").Value = 33
$ws.Range("E158").Value = 88
